$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells (outside edited block) for format-only paste

# Row 305: R.02.0022 -> O.01.0120
$ws.Range("H305").Value2 = "O.01.0120"
$ws.Range("I305").Value2 = "CHAPA PLASTIFICADO 20 MM - 2,20 X 1,10 M"
$ws.Range("J305").Value2 = "UN"
$ws.Range("O305").Value2 = "KALUTA"
$ws.Range("K305").Value2 = 14
$ws.Range("L305").Value2 = 160
$ws.Range("M305").Value2 = 2240
$ws.Range("F305").Value2 = 80785
$ws.Range("G305").Value2 = 45982
$ws.Range("N305").Value2 = "'00000000009484"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N305").PasteSpecial(-4122) | Out-Null

# Row 306: W.01.0060 -> O.01.0180
$ws.Range("H306").Value2 = "O.01.0180"
$ws.Range("I306").Value2 = "MADEIRA DE LEI"
$ws.Range("J306").Value2 = "M"
$ws.Range("O306").Value2 = "KALUTA"
$ws.Range("K306").Value2 = 35
$ws.Range("L306").Value2 = 33.6
$ws.Range("M306").Value2 = 1176
$ws.Range("F306").Value2 = 80785
$ws.Range("G306").Value2 = 45982
$ws.Range("N306").Value2 = "'00000000009484"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N306").PasteSpecial(-4122) | Out-Null

# Row 307: W.01.0006 -> O.01.0112
$ws.Range("H307").Value2 = "O.01.0112"
$ws.Range("I307").Value2 = "CHAPA PLASTIFICADO 10 MM - 2,20 X 1,10 M"
$ws.Range("J307").Value2 = "UN"
$ws.Range("O307").Value2 = "KALUTA"
$ws.Range("K307").Value2 = 40
$ws.Range("L307").Value2 = 82
$ws.Range("M307").Value2 = 3280
$ws.Range("F307").Value2 = 80785
$ws.Range("G307").Value2 = 45982
$ws.Range("N307").Value2 = "'00000000009484"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N307").PasteSpecial(-4122) | Out-Null

# Row 308: W.01.0008 -> O.01.0004
$ws.Range("H308").Value2 = "O.01.0004"
$ws.Range("I308").Value2 = "TÁBUA DE CEDRINHO - 1 X 4'' -"
$ws.Range("J308").Value2 = "M"
$ws.Range("O308").Value2 = "KALUTA"
$ws.Range("K308").Value2 = 120
$ws.Range("L308").Value2 = 10
$ws.Range("M308").Value2 = 1200
$ws.Range("F308").Value2 = 80785
$ws.Range("G308").Value2 = 45982
$ws.Range("N308").Value2 = "'00000000009484"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N308").PasteSpecial(-4122) | Out-Null

# Row 309: E.01.0037 -> O.01.0105
$ws.Range("H309").Value2 = "O.01.0105"
$ws.Range("I309").Value2 = "PONTALETE DE CEDRINHO - 3 X 3`" - 1ª IND"
$ws.Range("J309").Value2 = "M"
$ws.Range("O309").Value2 = "KALUTA"
$ws.Range("K309").Value2 = 180
$ws.Range("L309").Value2 = 14.4
$ws.Range("M309").Value2 = 2592
$ws.Range("F303").Copy() | Out-Null
$ws.Range("F309").PasteSpecial(-4122) | Out-Null
$ws.Range("F309").Value2 = 80785
$ws.Range("G303").Copy() | Out-Null
$ws.Range("G309").PasteSpecial(-4122) | Out-Null
$ws.Range("G309").Value2 = 45982
$ws.Range("N309").Value2 = "'00000000009484"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N309").PasteSpecial(-4122) | Out-Null

# Row 310: M.05.0200 -> O.01.0008
$ws.Range("H310").Value2 = "O.01.0008"
$ws.Range("I310").Value2 = "TÁBUA DE CEDRINHO - 1 X 12'' -"
$ws.Range("J310").Value2 = "M"
$ws.Range("O310").Value2 = "KALUTA"
$ws.Range("K310").Value2 = 80
$ws.Range("L310").Value2 = 30
$ws.Range("M310").Value2 = 2400
$ws.Range("F303").Copy() | Out-Null
$ws.Range("F310").PasteSpecial(-4122) | Out-Null
$ws.Range("F310").Value2 = 80785
$ws.Range("G303").Copy() | Out-Null
$ws.Range("G310").PasteSpecial(-4122) | Out-Null
$ws.Range("G310").Value2 = 45982
$ws.Range("N310").Value2 = "'00000000009484"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N310").PasteSpecial(-4122) | Out-Null

# Row 311: O.01.0120 -> R.02.0022
$ws.Range("H311").Value2 = "R.02.0022"
$ws.Range("I311").Value2 = "TINTA ACRILICA LATA DE 18L"
$ws.Range("J311").Value2 = "LAT"
$ws.Range("O311").Value2 = "REI DAS TINTAS"
$ws.Range("K311").Value2 = 2
$ws.Range("L311").Value2 = 530
$ws.Range("M311").Value2 = 1060
$ws.Range("F303").Copy() | Out-Null
$ws.Range("F311").PasteSpecial(-4122) | Out-Null
$ws.Range("F311").Value2 = 80760
$ws.Range("G303").Copy() | Out-Null
$ws.Range("G311").PasteSpecial(-4122) | Out-Null
$ws.Range("G311").Value2 = 45981
$ws.Range("N311").Value2 = "'00000000008674"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N311").PasteSpecial(-4122) | Out-Null

# Row 312: O.01.0180 -> W.01.0060
$ws.Range("H312").Value2 = "W.01.0060"
$ws.Range("I312").Value2 = "PREGO COM ROSCA SOBERBA  PARA TELHA   COM VEDAÇÃO"
$ws.Range("J312").Value2 = "KG"
$ws.Range("O312").Value2 = "SPW3"
$ws.Range("K312").Value2 = 2
$ws.Range("L312").Value2 = 42
$ws.Range("M312").Value2 = 84
$ws.Range("F303").Copy() | Out-Null
$ws.Range("F312").PasteSpecial(-4122) | Out-Null
$ws.Range("F312").Value2 = 80754
$ws.Range("G303").Copy() | Out-Null
$ws.Range("G312").PasteSpecial(-4122) | Out-Null
$ws.Range("G312").Value2 = 45980
$ws.Range("N312").Value2 = "'00000000002393"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N312").PasteSpecial(-4122) | Out-Null

# Row 313: O.01.0112 -> W.01.0006
$ws.Range("H313").Value2 = "W.01.0006"
$ws.Range("I313").Value2 = "PREGO COMUM C/ CABEÇA 17 X 27"
$ws.Range("J313").Value2 = "KG"
$ws.Range("O313").Value2 = "SPW3"
$ws.Range("K313").Value2 = 10
$ws.Range("L313").Value2 = 17
$ws.Range("M313").Value2 = 170
$ws.Range("F303").Copy() | Out-Null
$ws.Range("F313").PasteSpecial(-4122) | Out-Null
$ws.Range("F313").Value2 = 80754
$ws.Range("G303").Copy() | Out-Null
$ws.Range("G313").PasteSpecial(-4122) | Out-Null
$ws.Range("G313").Value2 = 45980
$ws.Range("N313").Value2 = "'00000000002393"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N313").PasteSpecial(-4122) | Out-Null

# Row 314: O.01.0004 -> W.01.0008
$ws.Range("H314").Value2 = "W.01.0008"
$ws.Range("I314").Value2 = "PREGO COMUM C/ CABEÇA 19 X 36"
$ws.Range("J314").Value2 = "KG"
$ws.Range("O314").Value2 = "SPW3"
$ws.Range("K314").Value2 = 3
$ws.Range("L314").Value2 = 22
$ws.Range("M314").Value2 = 66
$ws.Range("F303").Copy() | Out-Null
$ws.Range("F314").PasteSpecial(-4122) | Out-Null
$ws.Range("F314").Value2 = 80754
$ws.Range("G303").Copy() | Out-Null
$ws.Range("G314").PasteSpecial(-4122) | Out-Null
$ws.Range("G314").Value2 = 45980
$ws.Range("N314").Value2 = "'00000000002393"
$ws.Range("N303").Copy() | Out-Null
$ws.Range("N314").PasteSpecial(-4122) | Out-Null

# Row 315: O.01.0105 -> E.01.0037
$ws.Range("H315").Value2 = "E.01.0037"
$ws.Range("I315").Value2 = "LONA AZUL  REFORÇADA COM OLHAL"
$ws.Range("J315").Value2 = "UN"
$ws.Range("K315").Value2 = 1

# Row 316: O.01.0008 -> M.05.0200
$ws.Range("H316").Value2 = "M.05.0200"
$ws.Range("I316").Value2 = "TELHA METÁLICA TRAPEZOIDAL COM PINTURA  EM DUAS FACES"
$ws.Range("J316").Value2 = "M²"
$ws.Range("K316").Value2 = 40

$excel.CutCopyMode = 0